# Weekly update: insert a new price record as row 10, pushing the
# previously-existing rows 10-14 down to 11-15 (data itself is unchanged,
# only its row position shifts).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 10.
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with this week's record.
$ws.Cells.Item(10, 1).Value  = 6
$ws.Cells.Item(10, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(10, 3).Value  = "Metropolitana"
$ws.Cells.Item(10, 4).Value  = 44483
$ws.Cells.Item(10, 5).Value  = 13
$ws.Cells.Item(10, 6).Value  = 100112035
$ws.Cells.Item(10, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(10, 8).Value  = "Sin especificar"
$ws.Cells.Item(10, 9).Value  = "Primera"
$ws.Cells.Item(10, 10).Value = 220
$ws.Cells.Item(10, 11).Value = 18000
$ws.Cells.Item(10, 12).Value = 20000
$ws.Cells.Item(10, 13).Value = 18909
$ws.Cells.Item(10, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(10, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(10, 16).Value = 1261
$ws.Cells.Item(10, 17).Value = 15
$ws.Cells.Item(10, 18).Value = "Hortaliza"
